$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: Remove the "_GoBack" bookmark from its original location
# (end of paragraph 2). It will be re-added at the end of the document
# once the new paragraphs are in place.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# Step 2: Rewrite paragraph 2's text as a single run (the original had
# 3 runs split around a <w:proofErr> pair). Word's COM layer only
# rebuilds run structure when the text actually differs from what is
# already there, so first stamp a throw-away placeholder, then set the
# real text.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$r = $p2.Range
$r.End = $r.End - 1
$r.Text = "x"
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Text = "Toutes les pages du site ont comme modèle la page maître «HobbyCartes.master». Celui-ci contient l’entête, le menu et le pied de page qui sont identiques d’une page à l’autre."

# ---------------------------------------------------------------------
# Step 3: Append four more list paragraphs (same "Paragraphedeliste"
# list style) after paragraph 2.
# ---------------------------------------------------------------------
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Text = "Tous les onglets de la page membre sont des pages qui ont comme modèle la page maître «Membre.master» qui a lui-même comme modèle «HobbyCartes.master»."

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = "Les noms d’état acceptés pour une fiche sont : «impeccable, bonne, moyenne, passable et pietre». Il faut respecter la casse et les accents."

# Paragraph 5 needs three distinct runs with identical formatting.
# Word's save pass coalesces adjacent same-format runs within one
# paragraph, so build each chunk of text as its own paragraph first,
# then join them by deleting the paragraph mark between them -- that
# keeps the runs separate instead of merging them.
$p4.Range.InsertParagraphAfter()
$p5a = $d.Paragraphs(5)
$r5a = $p5a.Range
$r5a.End = $r5a.End - 1
$r5a.Text = "Les types de collections acceptés sont "

$p5a.Range.InsertParagraphAfter()
$p5b = $d.Paragraphs(6)
$r5b = $p5b.Range
$r5b.End = $r5b.End - 1
$r5b.Text = ": «"

$p5b.Range.InsertParagraphAfter()
$p5c = $d.Paragraphs(7)
$r5c = $p5c.Range
$r5c.End = $r5c.End - 1
$r5c.Text = "hockey, football, baseball et basketball». Il faut respecter la casse."

# Join paragraph 5a+5b (delete the paragraph mark ending 5a).
$p5aEnd = $d.Paragraphs(5).Range.End
$d.Range($p5aEnd - 1, $p5aEnd).Delete()
# Join the (now merged) 5a/5b with what used to be 5c.
$p5End = $d.Paragraphs(5).Range.End
$d.Range($p5End - 1, $p5End).Delete()

# ---------------------------------------------------------------------
# Step 4: Final list paragraph.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$r6.End = $r6.End - 1
$r6.Text = "La page «acceuil.aspx» se nomme «Default.aspx» dans le projet."

# ---------------------------------------------------------------------
# Step 5: Re-add the "_GoBack" bookmark at the end of the last
# paragraph (collapsed, right before its paragraph mark), matching
# where it now lives in the edited document.
#
# Word's COM layer here mis-handles Bookmarks.Add for a collapsed range
# that sits exactly at "end of paragraph text, just before the mark"
# (it snaps to the whole-paragraph range instead). Work around it by
# inserting a throw-away trailing character, anchoring the collapsed
# bookmark just before that character (now a safe mid-paragraph spot),
# then deleting the throw-away character.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
$tailRange = $p6.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter("Z")

$p6b = $d.Paragraphs(6)
$bmRange = $p6b.Range
$bmRange.MoveEnd(1, -1)
$bmRange.MoveEnd(1, -1)
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

$bm = $d.Bookmarks("_GoBack")
$zPos = $bm.End
$d.Range($zPos, $zPos + 1).Delete()

Write-Output "done"
